$d = $word.ActiveDocument

$d.Content.Find.Execute("713÷3=237, 2", $true, $false, $false, $false, $false, $true, 1, $false, "303÷3=101, 0", 2) | Out-Null
$d.Content.Find.Execute("574÷3=191, 1", $true, $false, $false, $false, $false, $true, 1, $false, "946÷4=236, 2", 2) | Out-Null
$d.Content.Find.Execute("136÷2=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "988÷3=329, 1", 2) | Out-Null
$d.Content.Find.Execute("374÷4=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "773÷4=193, 1", 2) | Out-Null
$d.Content.Find.Execute("455÷5=91, 0", $true, $false, $false, $false, $false, $true, 1, $false, "520÷3=173, 1", 2) | Out-Null
$d.Content.Find.Execute("699÷9=77, 6", $true, $false, $false, $false, $false, $true, 1, $false, "981÷2=490, 1", 2) | Out-Null
$d.Content.Find.Execute("179÷7=25, 4", $true, $false, $false, $false, $false, $true, 1, $false, "626÷5=125, 1", 2) | Out-Null
$d.Content.Find.Execute("154÷6=25, 4", $true, $false, $false, $false, $false, $true, 1, $false, "800÷2=400, 0", 2) | Out-Null
$d.Content.Find.Execute("867÷7=123, 6", $true, $false, $false, $false, $false, $true, 1, $false, "121÷7=17, 2", 2) | Out-Null
$d.Content.Find.Execute("853÷6=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "467÷2=233, 1", 2) | Out-Null
$d.Content.Find.Execute("532÷6=88, 4", $true, $false, $false, $false, $false, $true, 1, $false, "290÷7=41, 3", 2) | Out-Null
$d.Content.Find.Execute("432÷8=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "103÷7=14, 5", 2) | Out-Null
$d.Content.Find.Execute("958÷6=159, 4", $true, $false, $false, $false, $false, $true, 1, $false, "812÷3=270, 2", 2) | Out-Null
$d.Content.Find.Execute("736÷6=122, 4", $true, $false, $false, $false, $false, $true, 1, $false, "812÷9=90, 2", 2) | Out-Null
$d.Content.Find.Execute("995÷6=165, 5", $true, $false, $false, $false, $false, $true, 1, $false, "375÷5=75, 0", 2) | Out-Null
$d.Content.Find.Execute("606÷2=303, 0", $true, $false, $false, $false, $false, $true, 1, $false, "178÷4=44, 2", 2) | Out-Null
$d.Content.Find.Execute("722÷4=180, 2", $true, $false, $false, $false, $false, $true, 1, $false, "808÷4=202, 0", 2) | Out-Null
$d.Content.Find.Execute("335÷6=55, 5", $true, $false, $false, $false, $false, $true, 1, $false, "444÷9=49, 3", 2) | Out-Null
$d.Content.Find.Execute("135÷9=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "118÷5=23, 3", 2) | Out-Null
$d.Content.Find.Execute("717÷6=119, 3", $true, $false, $false, $false, $false, $true, 1, $false, "575÷9=63, 8", 2) | Out-Null
$d.Content.Find.Execute("173÷9=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "439÷5=87, 4", 2) | Out-Null
$d.Content.Find.Execute("898÷6=149, 4", $true, $false, $false, $false, $false, $true, 1, $false, "401÷8=50, 1", 2) | Out-Null
$d.Content.Find.Execute("935÷7=133, 4", $true, $false, $false, $false, $false, $true, 1, $false, "900÷6=150, 0", 2) | Out-Null
$d.Content.Find.Execute("619÷6=103, 1", $true, $false, $false, $false, $false, $true, 1, $false, "997÷8=124, 5", 2) | Out-Null
$d.Content.Find.Execute("347÷5=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "989÷5=197, 4", 2) | Out-Null
